$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.202003359794617
$ws.Range("B1").Value = 1.284418106079102
$ws.Range("C1").Value = 6.883109092712402
$ws.Range("D1").Value = 2.156526327133179
$ws.Range("E1").Value = 1.160950064659119
